# Fixed missing page number on energy bill page
#
# The section had a "different first page" header/footer set (w:titlePg)
# plus separate even-page headers/footers, left over from when this
# chapter was its own sub-document. Because of that, the first page of
# this section fell back to the (blank) "first page" footer instead of
# the normal default footer that prints the "${LE} - Page #" page
# number, so the first page of this chapter showed no page number.
#
# Fix: stop using different first-page / odd-even headers and footers so
# every page uses the single default footer (which already contains the
# PAGE field), and clear out the now-unused header/footer stories.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Clear the (empty) header content for all three header variants - this
# section should not show any header text.
$sec.Headers.Item(1).Range.Delete()   # default/primary header
$sec.Headers.Item(2).Range.Delete()   # first-page header
$sec.Headers.Item(3).Range.Delete()   # even-page header

# Clear the even-page and first-page footers; the default (primary)
# footer already holds the "${LE} - Page #" text and is left untouched
# so it keeps printing on every page.
$sec.Footers.Item(2).Range.Delete()   # first-page footer
$sec.Footers.Item(3).Range.Delete()   # even-page footer

# Stop treating the first page (and even pages) differently so the
# default footer's page number shows up on every page, including the
# first page of this section.
$sec.PageSetup.DifferentFirstPageHeaderFooter = $false
$sec.PageSetup.OddAndEvenPagesHeaderFooter = $false
